# Apply "edits to attribute tables" changes to the ColumnHeaders sheet.
#
# Summary of change:
#   - rename "lat" -> "latitude" (row 4, col A)
#   - rename "lon" -> "longitude" (row 5, col A)
#   - insert a new row (depth) between the lon/longitude row and the temp row,
#     shifting temp/sal/biosat/O2_Ar_ratio down by one row
#   - rename "O2_Ar_ratio" -> "O2_Ar_ratio_corrected" and update its
#     definition text (now the last row, row 10)
#   - update the active selection to A6:D6 (the newly inserted depth row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename lat -> latitude, lon -> longitude (values stay otherwise the same)
$ws.Range("A4").Value = "latitude"
$ws.Range("A5").Value = "longitude"

# Insert a new row 6 (pushes temp/sal/biosat/O2_Ar_ratio rows down by one)
$ws.Rows.Item(6).Insert()

# New row 6: depth attribute
$ws.Range("A6").Value = "depth"
$ws.Range("B6").Value = "Depth of sample below sea surface. URI http://vocab.nerc.ac.uk/collection/P09/current/DEPH/"
$ws.Range("C6").Value = "numeric"
$ws.Range("D6").Value = "meter"

# Last row (was row 9 "O2_Ar_ratio", now row 10) gets renamed + redefined
$ws.Range("A10").Value = "O2_Ar_ratio_corrected"
$ws.Range("B10").Value = "Oxygen-argon ratio of EIMS sample from underway corrected for air values"

# Update the selection to match the saved state
$ws.Range("A6:D6").Select() | Out-Null
